# ---------------------------------------------------------------------------
# Applies two changes described by the source diff:
#
#  1. Slide 5's table (B1- TYPES OF FINANCIAL DOCUMENTS) switches its table
#     style from the custom "Table_0" style to the built-in
#     "No Style, Table Grid" style ({CB60238C-2F9B-4ECC-9312-C650FA6D88E3}).
#
#  2. The presentation's theme ("Integral" / Red Violet) is swapped back to
#     the plain default "Office Theme" colour scheme (the deck's two theme
#     parts effectively trade places). The richest surface the PowerPoint
#     object model exposes for rewriting the live theme part in this host is
#     the 12-colour ThemeColorScheme, so the Office Theme palette is written
#     there, colour by colour, through a Slide's ThemeColorScheme (all
#     slides/layouts/master share the single theme part, so any slide works
#     as the entry point).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 ---------------------------------------------
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{CB60238C-2F9B-4ECC-9312-C650FA6D88E3}")
    }
}

# --- 2. Theme colour scheme -> default Office Theme palette ----------------
function Set-ThemeColor($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $scheme.Item($index).RGB = ($b * 65536) + ($g * 256) + $r
}

$officeTheme = @("000000", "FFFFFF", "44546A", "E7E6E6", "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47", "0563C1", "954F72")

$anySlide = $p.Slides.Item(1)
$colorScheme = $anySlide.ThemeColorScheme
for ($idx = 1; $idx -le $officeTheme.Count; $idx++) {
    Set-ThemeColor $colorScheme $idx $officeTheme[$idx - 1]
}
